$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$baseStyle = $ws.Range("B2").Style

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '58.288.22'
$ws.Range("D2").Style = $baseStyle
$ws.Range("E2").Value = '  +2.88%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.368.23'
$ws.Range("D3").Style = $baseStyle
$ws.Range("E3").Value = '  +1.39%  '

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.998'
$ws.Range("D4").Style = $baseStyle
$ws.Range("E4").Value = '  -0.25%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '544.91'
$ws.Range("D5").Style = $baseStyle
$ws.Range("E5").Value = '  +5.92%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '135.67'
$ws.Range("D6").Style = $baseStyle
$ws.Range("E6").Value = '  +1.75%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.998'
$ws.Range("D7").Style = $baseStyle
$ws.Range("E7").Value = '  -0.08%  '

# Row 8
$ws.Range("E8").Value = '  +0.84%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.363.70'
$ws.Range("D9").Style = $baseStyle

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.103'
$ws.Range("D10").Style = $baseStyle
$ws.Range("E10").Value = '  +1.74%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.41'
$ws.Range("D11").Style = $baseStyle
$ws.Range("E11").Value = '  +1.65%  '

# Row 12
$ws.Range("E12").Value = '  +0.66%  '

# Row 13
$ws.Range("E13").Value = '  +4.67%  '

# Row 14
$ws.Range("B14").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C14").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '2.787.80'
$ws.Range("D14").Style = $baseStyle
$ws.Range("E14").Value = '  +1.31%  '

# Row 15
$ws.Range("B15").Value = 'Avalanche'
$ws.Range("C15").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '23.73'
$ws.Range("D15").Style = $baseStyle
$ws.Range("E15").Value = '  +0.20%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '58.189.24'
$ws.Range("D16").Style = $baseStyle
$ws.Range("E16").Value = '  +2.76%  '

# Row 17
$ws.Range("E17").Value = '  +0.80%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.375.38'
$ws.Range("D18").Style = $baseStyle
$ws.Range("E18").Value = '  +1.54%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '339.37'
$ws.Range("D19").Style = $baseStyle
$ws.Range("E19").Value = '  +4.59%  '

# Row 20
$ws.Range("E20").Value = '  +1.17%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '4.24'
$ws.Range("D21").Style = $baseStyle
$ws.Range("E21").Value = '  +1.76%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.92'
$ws.Range("D22").Style = $baseStyle
$ws.Range("E22").Value = '  +4.09%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.999'
$ws.Range("D23").Style = $baseStyle
$ws.Range("E23").Value = '  +0.05%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '62.38'
$ws.Range("D24").Style = $baseStyle

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.169'
$ws.Range("D25").Style = $baseStyle
$ws.Range("E25").Value = '  +3.72%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '8.55'
$ws.Range("D26").Style = $baseStyle
$ws.Range("E26").Value = '  -0.98%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.998'
$ws.Range("D27").Style = $baseStyle
$ws.Range("E27").Value = '  +0.04%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.41'
$ws.Range("D28").Style = $baseStyle
$ws.Range("E28").Value = '  +7.92%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '174.19'
$ws.Range("D29").Style = $baseStyle
$ws.Range("E29").Value = '  +3.55%  '

# Row 30
$ws.Range("E30").Value = '  +5.38%  '

# Row 31
$ws.Range("E31").Value = '  +2.10%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '6.16'
$ws.Range("D32").Style = $baseStyle
$ws.Range("E32").Value = '  +0.76%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '18.58'
$ws.Range("D33").Style = $baseStyle
$ws.Range("E33").Value = '  +1.32%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.02'
$ws.Range("D34").Style = $baseStyle
$ws.Range("E34").Value = '  +15.07%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.999'
$ws.Range("D35").Style = $baseStyle
$ws.Range("E35").Value = '  -0.02%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.998'
$ws.Range("D36").Style = $baseStyle
$ws.Range("E36").Value = '  +0.07%  '

# Row 37
$ws.Range("E37").Value = '  -0.07%  '

# Row 38
$ws.Range("E38").Value = '  +3.42%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.62'
$ws.Range("D39").Style = $baseStyle
$ws.Range("E39").Value = '  +3.84%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '39.45'
$ws.Range("D40").Style = $baseStyle
$ws.Range("E40").Value = '  +2.63%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '150.76'
$ws.Range("D41").Style = $baseStyle
$ws.Range("E41").Value = '  +0.54%  '

# Row 42
$ws.Range("E42").Value = '  +0.97%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '3.64'
$ws.Range("D43").Style = $baseStyle
$ws.Range("E43").Value = '  +1.70%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '286.71'
$ws.Range("D44").Style = $baseStyle
$ws.Range("E44").Value = '  +2.53%  '

# Row 45
$ws.Range("B45").Value = 'InjectiveProtocol'
$ws.Range("C45").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '19.22'
$ws.Range("D45").Style = $baseStyle
$ws.Range("E45").Value = '  +5.65%  '

# Row 46
$ws.Range("B46").Value = 'Stellar'
$ws.Range("C46").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0929'
$ws.Range("D46").Style = $baseStyle
$ws.Range("E46").Value = '  +0.52%  '

# Row 47
$ws.Range("B47").Value = 'Hedera'
$ws.Range("C47").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0505'
$ws.Range("D47").Style = $baseStyle
$ws.Range("E47").Value = '  +1.52%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.561'
$ws.Range("D48").Style = $baseStyle
$ws.Range("E48").Value = '  +1.12%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0219'
$ws.Range("D49").Style = $baseStyle
$ws.Range("E49").Value = '  +1.77%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '17.57'
$ws.Range("D50").Style = $baseStyle
$ws.Range("E50").Value = '  +2.93%  '

# Row 51
$ws.Range("E51").Value = '  +0.67%  '
